$wb = $excel.ActiveWorkbook

# --- Sheet 1: LH-TC-REGISTERATION-Reviews ---
# Review comments for rows 7-9 (LH_TC_REGITERATION_REVIEW_007/008/009) are now
# closed: "Owner Status" (col I) and "Reviewer verification" (col J) move
# from "Open" to "Closed".
$ws1 = $wb.Worksheets.Item("LH-TC-REGISTERATION-Reviews")
$ws1.Range("I7").Value = "Closed"
$ws1.Range("J7").Value = "Closed"
$ws1.Range("I8").Value = "Closed"
$ws1.Range("J8").Value = "Closed"
$ws1.Range("I9").Value = "Closed"
$ws1.Range("J9").Value = "Closed"

# --- Sheet 2: Version History ---
# Add a new version history entry documenting this change.
$ws2 = $wb.Worksheets.Item("Version History")
$ws2.Range("A6").Value = "v2.1"
$ws2.Range("B6").Value = "Omar Sherif"
$ws2.Range("C6").Value = "review comments closed"
$ws2.Range("D6").Value = "13/5/2025"

# A6/B6 previously held empty, formatted cells; clear that formatting back to
# the workbook default (matches the unformatted A5/B5/D5 cells above them).
$ws2.Range("A6").Style = "Normal"
$ws2.Range("B6").Style = "Normal"

# Final selection / active sheet: cursor left on B8 of the reviews sheet, and
# the Version History tab is the one on top when the file is saved.
$ws1.Range("B8").Select()
$ws2.Range("E9").Select()
$ws2.Activate()
